# Applies the "Updated symbol list" data refresh (cryptos.xlsx, Sheet1).
# For each touched row, the B:E cells are re-written with the refreshed
# Coin / Link / Price / Volume(1h) values. D & E hold numeric/percent-looking
# text (e.g. "257.84", "0.47%") that Excel would otherwise auto-convert to a
# number; briefly forcing NumberFormat "@" (Text) before the write, then
# clearing formats back to the original (unstyled) state, keeps the cells as
# plain text without leaving any formatting behind - matching the workbook's
# original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '257.84'
$ws.Range("E2").Value = '0.47%'
$ws.Range("D2:E2").ClearFormats()

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '27.15'
$ws.Range("E3").Value = '-3.39%'
$ws.Range("D3:E3").ClearFormats()

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '4.873'
$ws.Range("E4").Value = '-7.83%'
$ws.Range("D4:E4").ClearFormats()

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05951'
$ws.Range("E5").Value = '2.44%'
$ws.Range("D5:E5").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.19%'
$ws.Range("E6").ClearFormats()

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8686'
$ws.Range("E7").Value = '-0.35%'
$ws.Range("D7:E7").ClearFormats()

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9592'
$ws.Range("E8").Value = '6.14%'
$ws.Range("D8:E8").ClearFormats()

$ws.Range("B9:E9").NumberFormat = "@"
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.0006041'
$ws.Range("E9").Value = '-94.26%'
$ws.Range("B9:E9").ClearFormats()

$ws.Range("B10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1413'
$ws.Range("E10").Value = '0.15%'
$ws.Range("B10:E10").ClearFormats()

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07182'
$ws.Range("E11").Value = '0.08%'
$ws.Range("D11:E11").ClearFormats()

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03151'
$ws.Range("E12").Value = '-0.43%'
$ws.Range("D12:E12").ClearFormats()

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09241'
$ws.Range("E13").Value = '0.09%'
$ws.Range("D13:E13").ClearFormats()

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001540'
$ws.Range("E14").Value = '0.41%'
$ws.Range("D14:E14").ClearFormats()

$ws.Range("B15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005984'
$ws.Range("E15").Value = '0.62%'
$ws.Range("B15:E15").ClearFormats()

$ws.Range("B16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.485'
$ws.Range("E16").Value = '-0.61%'
$ws.Range("B16:E16").ClearFormats()

$ws.Range("B17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '3.256'
$ws.Range("E17").Value = '0.97%'
$ws.Range("B17:E17").ClearFormats()

$ws.Range("B18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.205'
$ws.Range("E18").Value = '-2.92%'
$ws.Range("B18:E18").ClearFormats()

$ws.Range("B19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3146'
$ws.Range("E19").Value = '0.57%'
$ws.Range("B19:E19").ClearFormats()

$ws.Range("B20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '0.03540'
$ws.Range("E20").Value = '3.45%'
$ws.Range("B20:E20").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.49%'
$ws.Range("E21").ClearFormats()

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.546'
$ws.Range("E22").Value = '0.79%'
$ws.Range("D22:E22").ClearFormats()

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04280'
$ws.Range("E23").Value = '2.87%'
$ws.Range("D23:E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.53%'
$ws.Range("E24").ClearFormats()

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001223'
$ws.Range("E25").Value = '-0.39%'
$ws.Range("D25:E25").ClearFormats()

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004517'
$ws.Range("E26").Value = '-9.18%'
$ws.Range("D26:E26").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-22.95%'
$ws.Range("E28").ClearFormats()

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03834'
$ws.Range("E40").Value = '-0.81%'
$ws.Range("D40:E40").ClearFormats()

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006582'
$ws.Range("E41").Value = '13.72%'
$ws.Range("D41:E41").ClearFormats()

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1102'
$ws.Range("E42").Value = '0.47%'
$ws.Range("D42:E42").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-10.09%'
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.12%'
$ws.Range("E44").ClearFormats()

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005496'
$ws.Range("E45").Value = '4.09%'
$ws.Range("D45:E45").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.14%'
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '28.53%'
$ws.Range("E47").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-2.21%'
$ws.Range("E48").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.14%'
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.14%'
$ws.Range("E50").ClearFormats()
